# update xlsxTransform to vTransform
# Appends 20 new "special thanks" rows (A1:D240 -> A1:D260) to the single
# worksheet, introducing 23 new shared strings (+1 reused existing string).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Seed formatting for the new cells first (copy/paste-special formats
#    from existing, equivalently-styled cells) so the values we set next
#    land with the same cell styles (s="1" / s="2") the source file uses.
# ---------------------------------------------------------------------

# Column B on every new row (241-260) uses the same centred/font-2 style
# as the rest of the "name" column (e.g. B240).
$ws.Range("B240").Copy()
$ws.Range("B241:B260").PasteSpecial(-4122)

# Column C "comment" cells that need the font-2 style (rows 247 & 253) -
# copy from an existing fully-populated comment cell (C43).
$ws.Range("C43").Copy()
$ws.Range("C247").PasteSpecial(-4122)
$ws.Range("C253").PasteSpecial(-4122)

# Column D (row 258) also uses the font-2 style - copy from D43.
$ws.Range("D43").Copy()
$ws.Range("D258").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# (Note: C258 intentionally keeps the plain default style, so it is left
# untouched - new cells written via .Value default to that style already.)

# ---------------------------------------------------------------------
# 2) Write the new row values.
# ---------------------------------------------------------------------

$ws.Range("A241").Value = 2
$ws.Range("B241").Value = "飞鸟"

$ws.Range("A242").Value = 2
$ws.Range("B242").Value = "Caesar"

$ws.Range("A243").Value = 2
$ws.Range("B243").Value = "张凯"

$ws.Range("A244").Value = 2
$ws.Range("B244").Value = "宅男阿海"

$ws.Range("A245").Value = 2
$ws.Range("B245").Value = "程蝶衣"

$ws.Range("A246").Value = 2
$ws.Range("B246").Value = "八千"

$ws.Range("A247").Value = 1
$ws.Range("B247").Value = "b 安贝慧"
$ws.Range("C247").Value = "  "

$ws.Range("A248").Value = 2
$ws.Range("B248").Value = "Alex"

$ws.Range("A249").Value = 2
$ws.Range("B249").Value = "呼噜呼噜"

$ws.Range("A250").Value = 2
$ws.Range("B250").Value = "解无明"

$ws.Range("A251").Value = 2
$ws.Range("B251").Value = "千鹤"

$ws.Range("A252").Value = 2
$ws.Range("B252").Value = "zeqing"

$ws.Range("A253").Value = 1
$ws.Range("B253").Value = "Si_X:"
$ws.Range("C253").Value = "社畜一枚，玩着你们的游戏睡着了，我做了一个梦，这个梦真的很美，谢谢"

$ws.Range("A254").Value = 2
$ws.Range("B254").Value = "不喜欢甜的"

$ws.Range("A255").Value = 2
$ws.Range("B255").Value = "晓危"

$ws.Range("A256").Value = 2
$ws.Range("B256").Value = "毛虫哥哥"

$ws.Range("A257").Value = 2
$ws.Range("B257").Value = "佚名"

$ws.Range("A258").Value = 1
$ws.Range("B258").Value = "听风忆雪"
$ws.Range("C258").Value = "搞死我了"
$ws.Range("D258").Value = "#4169E1"

$ws.Range("A259").Value = 2
$ws.Range("B259").Value = "老郑"

$ws.Range("A260").Value = 2
$ws.Range("B260").Value = "太帅很苦恼"

# ---------------------------------------------------------------------
# 3) Update the view state to match the edited workbook (scrolled down
#    to the newly-added rows, with I233 selected).
# ---------------------------------------------------------------------

$aw = $excel.ActiveWindow
$aw.ScrollRow = 223
$aw.ScrollColumn = 1
[void]$ws.Range("I233").Select()
